# Apply weekly update: insert two new price rows for "Comercializadora del
# Agro de Limarí - Uva" at the top of the data block (row 204), pushing the
# existing rows down by two. Fill the two new rows with the latest week's
# price data (Crimpson Seedless and Red Globe varieties).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 204 (existing rows 204+ shift down to 206+)
$ws.Rows.Item(204).Insert()
$ws.Rows.Item(204).Insert()

# --- Row 204: Crimpson Seedless -------------------------------------------
$ws.Range("A204").Value = 2
$ws.Range("B204").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C204").Value = "Coquimbo"
$ws.Range("D204").Value = 45077
$ws.Range("E204").Value = 4
$ws.Range("F204").Value = "Fruta"
$ws.Range("G204").Value = 100109
$ws.Range("H204").Value = "Uva"
$ws.Range("I204").Value = 100109001
$ws.Range("J204").Value = "Uva"
$ws.Range("K204").Value = "Crimpson Seedless"
$ws.Range("L204").Value = "Primera"
$ws.Range("M204").Value = 400
$ws.Range("N204").Value = 11000
$ws.Range("O204").Value = 12000
$ws.Range("P204").Value = 11500
$ws.Range("Q204").Value = "$/bandeja 18 kilos"
$ws.Range("R204").Value = "Provincia de Limarí"
$ws.Range("S204").Value = 639
$ws.Range("T204").Value = 18

# --- Row 205: Red Globe -----------------------------------------------------
$ws.Range("A205").Value = 2
$ws.Range("B205").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C205").Value = "Coquimbo"
$ws.Range("D205").Value = 45077
$ws.Range("E205").Value = 4
$ws.Range("F205").Value = "Fruta"
$ws.Range("G205").Value = 100109
$ws.Range("H205").Value = "Uva"
$ws.Range("I205").Value = 100109001
$ws.Range("J205").Value = "Uva"
$ws.Range("K205").Value = "Red Globe"
$ws.Range("L205").Value = "Primera"
$ws.Range("M205").Value = 600
$ws.Range("N205").Value = 11000
$ws.Range("O205").Value = 12000
$ws.Range("P205").Value = 11500
$ws.Range("Q205").Value = "$/bandeja 18 kilos"
$ws.Range("R205").Value = "Provincia de Limarí"
$ws.Range("S205").Value = 639
$ws.Range("T205").Value = 18

# Ensure the date columns keep the same date number format as the rest of
# column D.
$ws.Range("D204").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D205").NumberFormat = "YYYY-MM-DD HH:MM:SS"
